$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style_2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.858.84'
$ws.Range("D2").Style = $style_2
$ws.Range("E2").Value = '  -2.42%  '

$style_3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.580.27'
$ws.Range("D3").Style = $style_3
$ws.Range("E3").Value = '  -0.13%  '

$ws.Range("E4").Value = '  +0.06%  '

$style_5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.65'
$ws.Range("D5").Style = $style_5
$ws.Range("E5").Value = '  +0.11%  '

$style_6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.01'
$ws.Range("D6").Style = $style_6
$ws.Range("E6").Value = '  +0.25%  '

$style_7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.576'
$ws.Range("D7").Style = $style_7
$ws.Range("E7").Value = '  -1.14%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  -2.83%  '

$style_10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.29'
$ws.Range("D10").Style = $style_10

$style_11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0814'
$ws.Range("D11").Style = $style_11
$ws.Range("E11").Value = '  -0.72%  '

$style_12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.65'
$ws.Range("D12").Style = $style_12
$ws.Range("E12").Value = '  -2.45%  '

$ws.Range("E13").Value = '  +6.12%  '

$style_14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.483.91'
$ws.Range("D14").Style = $style_14
$ws.Range("E14").Value = '  -3.68%  '

$style_15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.888'
$ws.Range("D15").Style = $style_15
$ws.Range("E15").Value = '  -0.74%  '

$style_16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.39'
$ws.Range("D16").Style = $style_16
$ws.Range("E16").Value = '  -0.24%  '

$style_17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '42.938.16'
$ws.Range("D17").Style = $style_17
$ws.Range("E17").Value = '  -2.12%  '

$ws.Range("B18").Value = 'InternetComputer(DFINITY)'
$ws.Range("C18").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$style_18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.93'
$ws.Range("D18").Style = $style_18
$ws.Range("E18").Value = '  +3.09%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$style_19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0997'
$ws.Range("D19").Style = $style_19
$ws.Range("E19").Value = '  +1.00%  '

$style_20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.65'
$ws.Range("D20").Style = $style_20
$ws.Range("E20").Value = '  -1.23%  '

$style_21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.09'
$ws.Range("D21").Style = $style_21
$ws.Range("E21").Value = '  -2.17%  '

$style_22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '254.77'
$ws.Range("D22").Style = $style_22
$ws.Range("E22").Value = '  -4.55%  '

$style_23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.95'
$ws.Range("D23").Style = $style_23
$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("E24").Value = '  -5.23%  '

$style_25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '28.93'
$ws.Range("D25").Style = $style_25
$ws.Range("E25").Value = '  -1.96%  '

$ws.Range("E26").Value = '  -0.19%  '

$style_27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.30'
$ws.Range("D27").Style = $style_27
$ws.Range("E27").Value = '  +0.15%  '

$style_28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.70'
$ws.Range("D28").Style = $style_28
$ws.Range("E28").Value = '  -1.83%  '

$ws.Range("E29").Value = '  -2.66%  '

$style_30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.04'
$ws.Range("D30").Style = $style_30
$ws.Range("E30").Value = '  -3.32%  '

$style_31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.04'
$ws.Range("D31").Style = $style_31
$ws.Range("E31").Value = '  +1.37%  '

$style_32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.17'
$ws.Range("D32").Style = $style_32
$ws.Range("E32").Value = '  -1.98%  '

$style_33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.42'
$ws.Range("D33").Style = $style_33
$ws.Range("E33").Value = '  -5.13%  '

$ws.Range("E34").Value = '  -1.79%  '

$style_35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0804'
$ws.Range("D35").Style = $style_35
$ws.Range("E35").Value = '  -1.60%  '

$style_36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.20'
$ws.Range("D36").Style = $style_36
$ws.Range("E36").Value = '  +7.70%  '

$ws.Range("E37").Value = '  -3.36%  '

$ws.Range("E38").Value = '  -0.80%  '

$style_39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.22'
$ws.Range("D39").Style = $style_39
$ws.Range("E39").Value = '  -3.70%  '

$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$style_40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.43'
$ws.Range("D40").Style = $style_40
$ws.Range("E40").Value = '  -4.55%  '

$ws.Range("E41").Value = '  -1.92%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$style_42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.88'
$ws.Range("D42").Style = $style_42
$ws.Range("E42").Value = '  -0.34%  '

$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$style_43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.07'
$ws.Range("D43").Style = $style_43
$ws.Range("E43").Value = '  +26.01%  '

$style_44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.075.15'
$ws.Range("D44").Style = $style_44
$ws.Range("E44").Value = '  +1.60%  '

$style_45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("D45").Style = $style_45

$style_46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.21'
$ws.Range("D46").Style = $style_46
$ws.Range("E46").Value = '  +0.27%  '

$style_47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.54'
$ws.Range("D47").Style = $style_47
$ws.Range("E47").Value = '  -3.02%  '

$style_48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '76.79'
$ws.Range("D48").Style = $style_48
$ws.Range("E48").Value = '  +10.06%  '

$style_49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '106.48'
$ws.Range("D49").Style = $style_49
$ws.Range("E49").Value = '  +0.59%  '

$style_50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.825.48'
$ws.Range("D50").Style = $style_50
$ws.Range("E50").Value = '  -0.28%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$style_51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.192'
$ws.Range("D51").Style = $style_51
$ws.Range("E51").Value = '  -0.41%  '
